$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing comment in F4 with additional text (extends the shared string
# used by F4 - the "country scale" comment on the Bhattari & Hammig row).
$ws.Range("F4").Value = "Doing this at the country scale might mean less variation? Also much longer time series. Also, no diagnostics apart from R2, and no statement of variation for fixed or random effects"

# Insert two new rows (5 and 6) before the old "Basse et al" row (previously
# row 6), pushing it down to row 8 and leaving row 7 empty.
$ws.Rows("5:6").Insert()

# Row 5: Bonilla-Bedoya et al entry
$ws.Range("A5").Value = "Bonilla-Bedoya et al"
$ws.Range("B5").Value = 2018
$ws.Range("C5").Value = "Socioecological system and potential deforestatino in Western Amazon forest landscapes"
$ws.Range("D5").Value = "Modelling potential land use change. Detemined whether there were relationships between vulnerability to forest loss and the management policies"
$ws.Range("E5").Value = "Uses biophysical and socioeconomic variables. Use maximum entropy model."
$ws.Range("F5").Value = "I think quite nice, although I need to read up a bit on entropy models. They cite Souza and De Marco 2014 who go into it."

# Row 6: Souza & De Marco entry (only Authors/Year/Title filled in so far)
$ws.Range("A6").Value = "Souza & De Marco"
$ws.Range("B6").Value = 2014
$ws.Range("C6").Value = "The use of species distribution models to predict the spatial distribution of deforestation in the western Brazilian Amazon"

# Match row heights of the two new rows to the wrapped-text content (57.6pt),
# matching the other wrapped rows in this table.
$ws.Range("A5:F5").RowHeight = 57.6
$ws.Range("A6:F6").RowHeight = 57.6

# Update the selected cell to D6, matching where editing left off.
$ws.Range("D6").Select()
